# edit.ps1 - Apply the "Updated cryptos list" GitHub Actions refresh to Sheet1.
# Updates the Price (column D) and Volume(1h) (column E) columns for the
# rows whose values changed, matching the target OOXML diff exactly.
#
# Certain new Price values look like plain decimal numbers (e.g. "8.08").
# Assigning those directly to Range.Value would cause Excel to auto-convert
# them into numeric cells (and even introduce floating point artifacts like
# 8.0800000000000001), whereas the source data stores them as literal text
# strings. To preserve the original "text number" representation we
# temporarily force those cells to text format before assigning, then
# restore their style afterwards so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: force text number-format on D-column cells whose new value would otherwise be auto-parsed as a number
$textFormatCells = @("D5", "D6", "D10", "D15", "D19", "D21", "D22", "D23", "D24", "D29", "D34", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D46", "D47", "D48", "D49", "D50")
foreach ($cell in $textFormatCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Step 2: write the new cell values
$ws.Range("D2").Value = '66.258.52'
$ws.Range("E2").Value = '  -0.74%  '
$ws.Range("D3").Value = '3.495.15'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D5").Value = '604.97'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").Value = '144.48'
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("D7").Value = '3.493.03'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("D10").Value = '8.08'
$ws.Range("E10").Value = '  +2.06%  '
$ws.Range("E11").Value = '  -4.18%  '
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").Value = '4.083.95'
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("E14").Value = '  -3.88%  '
$ws.Range("D15").Value = '30.38'
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("D16").Value = '3.489.36'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '66.277.40'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = '10.76'
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("E20").Value = '  -3.06%  '
$ws.Range("D21").Value = '14.87'
$ws.Range("E21").Value = '  -2.73%  '
$ws.Range("D22").Value = '426.59'
$ws.Range("E22").Value = '  -1.53%  '
$ws.Range("D23").Value = '0.594'
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("D24").Value = '77.98'
$ws.Range("D25").Value = '3.627.45'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("E28").Value = '  -5.43%  '
$ws.Range("D29").Value = '7.92'
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("E30").Value = '  -0.59%  '
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  -8.32%  '
$ws.Range("D34").Value = '25.08'
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("D35").Value = '3.479.27'
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").Value = '1.74'
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("D38").Value = '5.63'
$ws.Range("E38").Value = '  -4.57%  '
$ws.Range("D39").Value = '7.74'
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '169.99'
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("E42").Value = '  -3.46%  '
$ws.Range("D43").Value = '5.17'
$ws.Range("E43").Value = '  -4.44%  '
$ws.Range("D44").Value = '0.881'
$ws.Range("E44").Value = '  -1.59%  '
$ws.Range("E45").Value = '  -9.12%  '
$ws.Range("D46").Value = '45.44'
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("D47").Value = '25.97'
$ws.Range("E47").Value = '  -7.89%  '
$ws.Range("D48").Value = '1.21'
$ws.Range("E48").Value = '  -6.51%  '
$ws.Range("D49").Value = '2.43'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").Value = '7.15'
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("E51").Value = '  -2.78%  '

# Step 3: clear the temporary style so no stray formatting remains on the text cells
foreach ($cell in $textFormatCells) {
    $ws.Range($cell).Style = "Normal"
}
